$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows above the current row 239 (existing rows 239-299
# shift down to 241-301, carrying their formatting/styles with them).
$ws.Rows.Item(239).Resize(2).Insert()

# --- New row 239 ---
$ws.Range("A239").Value = 3
$ws.Range("B239").Value = "Femacal de La Calera"
$ws.Range("C239").Value = "Coquimbo"
$ws.Range("D239").Value = 44551
$ws.Range("E239").Value = 5
$ws.Range("F239").Value = 100112032
$ws.Range("G239").Value = "Zapallo italiano"
$ws.Range("H239").Value = "Sin especificar"
$ws.Range("I239").Value = "Primera"
$ws.Range("J239").Value = 185
$ws.Range("K239").Value = 4000
$ws.Range("L239").Value = 4500
$ws.Range("M239").Value = 4243
$ws.Range("N239").Value = "$/caja 36 unidades"
$ws.Range("O239").Value = "Provincia de Quillota"
$ws.Range("P239").Value = 118
$ws.Range("Q239").Value = 36
$ws.Range("R239").Value = "Hortaliza"

# --- New row 240 ---
$ws.Range("A240").Value = 3
$ws.Range("B240").Value = "Femacal de La Calera"
$ws.Range("C240").Value = "Coquimbo"
$ws.Range("D240").Value = 44551
$ws.Range("E240").Value = 5
$ws.Range("F240").Value = 100112032
$ws.Range("G240").Value = "Zapallo italiano"
$ws.Range("H240").Value = "Sin especificar"
$ws.Range("I240").Value = "Primera"
$ws.Range("J240").Value = 205
$ws.Range("K240").Value = 8000
$ws.Range("L240").Value = 8500
$ws.Range("M240").Value = 8273
$ws.Range("N240").Value = "$/caja 70 unidades"
$ws.Range("O240").Value = "Provincia de Quillota"
$ws.Range("P240").Value = 118
$ws.Range("Q240").Value = 70
$ws.Range("R240").Value = "Hortaliza"
